{"js": "// Replace each old formula/date text with its new counterpart.\n// Every \"old\" string is unique within the document, so a plain\n// (non-wildcard, case-sensitive) search safely targets exactly one run.\nconst replacements = [\n  [\"2024-05-16 Thursday\", \"2024-05-17 Friday\"],\n  [\"34+40=74\", \"6+93=99\"],\n  [\"43-24=19\", \"61-55=6\"],\n  [\"16+20=36\", \"91+7=98\"],\n  [\"86-80=6\", \"30+15=45\"],\n  [\"95-46=49\", \"30+58=88\"],\n  [\"8+43=51\", \"51-20=31\"],\n  [\"49+44=93\", \"3+15=18\"],\n  [\"13+27=40\", \"82-59=23\"],\n  [\"72+23=95\", \"90-32=58\"],\n  [\"58+4=62\", \"17+0=17\"],\n  [\"94+5=99\", \"74-17=57\"],\n  [\"12+21=33\", \"55-18=37\"],\n  [\"46-29=17\", \"98-71=27\"],\n  [\"73-57=16\", \"42+12=54\"],\n  [\"6+62=68\", \"60+29=89\"],\n  [\"44+1=45\", \"0+41=41\"],\n  [\"63-22=41\", \"69-18=51\"],\n  [\"91-66=25\", \"67-64=3\"],\n  [\"49-28=21\", \"10+19=29\"],\n  [\"81-73=8\", \"54-46=8\"],\n  [\"41+19=60\", \"22+50=72\"],\n  [\"75+11=86\", \"74-34=40\"],\n  [\"35+32=67\", \"97-48=49\"],\n  [\"86-60=26\", \"90-50=40\"],\n  [\"46+13=59\", \"10-9=1\"],\n  [\"32+6=38\", \"55+43=98\"],\n  [\"69-9=60\", \"63+8=71\"],\n  [\"3+90=93\", \"58-10=48\"],\n  [\"83-15=68\", \"54-36=18\"],\n  [\"62+28=90\", \"81-18=63\"],\n  [\"56-44=12\", \"4-1=3\"],\n  [\"42-21=21\", \"75-69=6\"],\n  [\"25+11=36\", \"44-29=15\"],\n  [\"19-6=13\", \"69-16=53\"],\n  [\"49+48=97\", \"23+46=69\"],\n  [\"68+2=70\", \"56-51=5\"],\n  [\"24+33=57\", \"9+80=89\"],\n  [\"11+86=97\", \"72+12=84\"],\n  [\"47-20=27\", \"9+38=47\"],\n  [\"70+19=89\", \"61+24=85\"],\n  [\"3+45=48\", \"44+44=88\"],\n  [\"59-21=38\", \"86-7=79\"],\n  [\"64-51=13\", \"29-8=21\"],\n  [\"64-1=63\", \"96-31=65\"],\n  [\"44+3=47\", \"22+61=83\"],\n  [\"34+16=50\", \"63-18=45\"],\n  [\"48-37=11\", \"36+30=66\"],\n  [\"6+16=22\", \"13+58=71\"],\n  [\"95-35=60\", \"88-7=81\"],\n  [\"83-7=76\", \"65-0=65\"],\n  [\"28+65=93\", \"49-43=6\"],\n  [\"84-76=8\", \"47-40=7\"],\n  [\"11+62=73\", \"54+22=76\"],\n  [\"78-59=19\", \"61-39=22\"],\n  [\"42+44=86\", \"34+2=36\"],\n  [\"76-20=56\", \"54+16=70\"],\n  [\"61+35=96\", \"3+42=45\"],\n  [\"76-62=14\", \"29-9=20\"],\n  [\"22+76=98\", \"87-78=9\"],\n  [\"47-38=9\", \"49-10=39\"],\n  [\"26+8=34\", \"59-25=34\"],\n  [\"55-14=41\", \"29+20=49\"],\n  [\"5+79=84\", \"21+31=52\"],\n  [\"11+27=38\", \"37-1=36\"],\n  [\"10+0=10\", \"85-67=18\"],\n  [\"22+67=89\", \"83-2=81\"],\n  [\"52-8=44\", \"15+43=58\"],\n  [\"61+11=72\", \"61-24=37\"],\n  [\"54+7=61\", \"20+65=85\"],\n  [\"42+26=68\", \"1+52=53\"],\n  [\"31+60=91\", \"13+48=61\"],\n  [\"93-13=80\", \"54-9=45\"],\n  [\"25-18=7\", \"81-60=21\"],\n  [\"49-37=12\", \"20+70=90\"],\n  [\"42-17=25\", \"78-47=31\"],\n  [\"74+4=78\", \"70+3=73\"],\n  [\"78-24=54\", \"14+36=50\"],\n  [\"28+17=45\", \"36-21=15\"],\n  [\"69+8=77\", \"2+1=3\"],\n  [\"6+23=29\", \"28+67=95\"],\n  [\"59+8=67\", \"52+42=94\"],\n  [\"56+43=99\", \"25+74=99\"],\n  [\"1+79=80\", \"4+69=73\"],\n  [\"13+42=55\", \"49+16=65\"],\n  [\"77+13=90\", \"59-10=49\"],\n  [\"99-5=94\", \"67+29=96\"],\n  [\"29+52=81\", \"3+23=26\"],\n  [\"93-16=77\", \"18+26=44\"],\n  [\"94-87=7\", \"58-58=0\"],\n  [\"53+25=78\", \"17+15=32\"],\n  [\"31-1=30\", \"60-11=49\"],\n  [\"48-20=28\", \"54+23=77\"],\n  [\"53+14=67\", \"76+7=83\"],\n  [\"34-33=1\", \"89-56=33\"],\n  [\"29+19=48\", \"5+55=60\"],\n  [\"10+73=83\", \"98-2=96\"],\n  [\"57-19=38\", \"68-7=61\"],\n  [\"39+37=76\", \"19+78=97\"],\n  [\"52-1=51\", \"76-50=26\"],\n  [\"47-22=25\", \"50-43=7\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    // insertText with Replace keeps the existing run formatting\n    // (font, size, etc.) and only swaps the text content.\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn \"replaced \" + replacements.length + \" items\";\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace each old formula/date text with its new counterpart.\n# Every \"old\" string is unique within the document, so a plain\n# (non-wildcard, case-sensitive) Find/Replace safely targets exactly\n# one run and leaves its formatting (font, size, etc.) untouched.\n$replacements = @(\n    @('2024-05-16 Thursday', '2024-05-17 Friday'),\n    @('34+40=74', '6+93=99'),\n    @('43-24=19', '61-55=6'),\n    @('16+20=36', '91+7=98'),\n    @('86-80=6', '30+15=45'),\n    @('95-46=49', '30+58=88'),\n    @('8+43=51', '51-20=31'),\n    @('49+44=93', '3+15=18'),\n    @('13+27=40', '82-59=23'),\n    @('72+23=95', '90-32=58'),\n    @('58+4=62', '17+0=17'),\n    @('94+5=99', '74-17=57'),\n    @('12+21=33', '55-18=37'),\n    @('46-29=17', '98-71=27'),\n    @('73-57=16', '42+12=54'),\n    @('6+62=68', '60+29=89'),\n    @('44+1=45', '0+41=41'),\n    @('63-22=41', '69-18=51'),\n    @('91-66=25', '67-64=3'),\n    @('49-28=21', '10+19=29'),\n    @('81-73=8', '54-46=8'),\n    @('41+19=60', '22+50=72'),\n    @('75+11=86', '74-34=40'),\n    @('35+32=67', '97-48=49'),\n    @('86-60=26', '90-50=40'),\n    @('46+13=59', '10-9=1'),\n    @('32+6=38', '55+43=98'),\n    @('69-9=60', '63+8=71'),\n    @('3+90=93', '58-10=48'),\n    @('83-15=68', '54-36=18'),\n    @('62+28=90', '81-18=63'),\n    @('56-44=12', '4-1=3'),\n    @('42-21=21', '75-69=6'),\n    @('25+11=36', '44-29=15'),\n    @('19-6=13', '69-16=53'),\n    @('49+48=97', '23+46=69'),\n    @('68+2=70', '56-51=5'),\n    @('24+33=57', '9+80=89'),\n    @('11+86=97', '72+12=84'),\n    @('47-20=27', '9+38=47'),\n    @('70+19=89', '61+24=85'),\n    @('3+45=48', '44+44=88'),\n    @('59-21=38', '86-7=79'),\n    @('64-51=13', '29-8=21'),\n    @('64-1=63', '96-31=65'),\n    @('44+3=47', '22+61=83'),\n    @('34+16=50', '63-18=45'),\n    @('48-37=11', '36+30=66'),\n    @('6+16=22', '13+58=71'),\n    @('95-35=60', '88-7=81'),\n    @('83-7=76', '65-0=65'),\n    @('28+65=93', '49-43=6'),\n    @('84-76=8', '47-40=7'),\n    @('11+62=73', '54+22=76'),\n    @('78-59=19', '61-39=22'),\n    @('42+44=86', '34+2=36'),\n    @('76-20=56', '54+16=70'),\n    @('61+35=96', '3+42=45'),\n    @('76-62=14', '29-9=20'),\n    @('22+76=98', '87-78=9'),\n    @('47-38=9', '49-10=39'),\n    @('26+8=34', '59-25=34'),\n    @('55-14=41', '29+20=49'),\n    @('5+79=84', '21+31=52'),\n    @('11+27=38', '37-1=36'),\n    @('10+0=10', '85-67=18'),\n    @('22+67=89', '83-2=81'),\n    @('52-8=44', '15+43=58'),\n    @('61+11=72', '61-24=37'),\n    @('54+7=61', '20+65=85'),\n    @('42+26=68', '1+52=53'),\n    @('31+60=91', '13+48=61'),\n    @('93-13=80', '54-9=45'),\n    @('25-18=7', '81-60=21'),\n    @('49-37=12', '20+70=90'),\n    @('42-17=25', '78-47=31'),\n    @('74+4=78', '70+3=73'),\n    @('78-24=54', '14+36=50'),\n    @('28+17=45', '36-21=15'),\n    @('69+8=77', '2+1=3'),\n    @('6+23=29', '28+67=95'),\n    @('59+8=67', '52+42=94'),\n    @('56+43=99', '25+74=99'),\n    @('1+79=80', '4+69=73'),\n    @('13+42=55', '49+16=65'),\n    @('77+13=90', '59-10=49'),\n    @('99-5=94', '67+29=96'),\n    @('29+52=81', '3+23=26'),\n    @('93-16=77', '18+26=44'),\n    @('94-87=7', '58-58=0'),\n    @('53+25=78', '17+15=32'),\n    @('31-1=30', '60-11=49'),\n    @('48-20=28', '54+23=77'),\n    @('53+14=67', '76+7=83'),\n    @('34-33=1', '89-56=33'),\n    @('29+19=48', '5+55=60'),\n    @('10+73=83', '98-2=96'),\n    @('57-19=38', '68-7=61'),\n    @('39+37=76', '19+78=97'),\n    @('52-1=51', '76-50=26'),\n    @('47-22=25', '50-43=7'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}\n\nWrite-Output \"replaced $($replacements.Count) items\"\n"}
